# Fruta / hortaliza, semanal
# Insert 3 new weekly records at the top of the Chirimoya - Macroferia Regional
# de Talca data block (rows 51-53), pushing the previously existing rows
# (old 51..84) down to (new 54..87).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 51..84 down by 3 rows, inserting 3 blank rows at 51.
$ws.Rows.Item(51).Resize(3).Insert()

# Common (unchanged across this data block) column values.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "$/bandeja 10 kilos"
$origen      = "Provincia de Limarí"
$kgPorUnidad = 10

# New weekly rows, with date 2022-09-09 (serial 44813).
$newRows = @(
    @{ Row = 51; Fecha = 44813; Calidad = "Especial"; Volumen = 100; Precio = 30000; PrecioKg = 3000 },
    @{ Row = 52; Fecha = 44813; Calidad = "Primera";  Volumen = 130; Precio = 25000; PrecioKg = 2500 },
    @{ Row = 53; Fecha = 44813; Calidad = "Segunda";  Volumen = 60;  Precio = 20000; PrecioKg = 2000 }
)

foreach ($r in $newRows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $mercadoId
    $ws.Cells.Item($row, 2).Value = $mercado
    $ws.Cells.Item($row, 3).Value = $region
    $ws.Cells.Item($row, 4).Value2 = $r.Fecha
    $ws.Cells.Item($row, 5).Value = $codreg
    $ws.Cells.Item($row, 6).Value = $tipo
    $ws.Cells.Item($row, 7).Value = $productoId
    $ws.Cells.Item($row, 8).Value = $producto
    $ws.Cells.Item($row, 9).Value = $categoriaId
    $ws.Cells.Item($row, 10).Value = $categoria
    $ws.Cells.Item($row, 11).Value = $variedad
    $ws.Cells.Item($row, 12).Value = $r.Calidad
    $ws.Cells.Item($row, 13).Value = $r.Volumen
    $ws.Cells.Item($row, 14).Value = $r.Precio
    $ws.Cells.Item($row, 15).Value = $r.Precio
    $ws.Cells.Item($row, 16).Value = $r.Precio
    $ws.Cells.Item($row, 17).Value = $unidad
    $ws.Cells.Item($row, 18).Value = $origen
    $ws.Cells.Item($row, 19).Value = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value = $kgPorUnidad
}
